$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue "D2" "43.956.42"
Set-TextValue "E2" "  +0.32%  "
Set-TextValue "D3" "2.249.29"
Set-TextValue "E3" "  -1.39%  "
Set-TextValue "E4" "  -0.04%  "
Set-TextValue "D5" "232.44"
Set-TextValue "E5" "  +0.67%  "
Set-TextValue "D6" "0.632"
Set-TextValue "E6" "  +0.72%  "
Set-TextValue "D7" "63.03"
Set-TextValue "E7" "  -2.33%  "
Set-TextValue "E8" "  -0.05%  "
Set-TextValue "E9" "  +5.58%  "
Set-TextValue "D10" "0.0981"
Set-TextValue "E10" "  +3.03%  "
Set-TextValue "D11" "57.33"
Set-TextValue "E11" "  -0.58%  "
Set-TextValue "D12" "26.37"
Set-TextValue "E12" "  -0.99%  "
Set-TextValue "D13" "0.106"
Set-TextValue "E13" "  +1.59%  "
Set-TextValue "D14" "2.582.14"
Set-TextValue "E14" "  -1.37%  "
Set-TextValue "D15" "15.48"
Set-TextValue "E15" "  -1.46%  "
Set-TextValue "E16" "  +3.11%  "
Set-TextValue "E17" "  +1.76%  "
Set-TextValue "D18" "2.252.23"
Set-TextValue "E18" "  -1.25%  "
Set-TextValue "D19" "43.853.02"
Set-TextValue "E19" "  +0.40%  "
Set-TextValue "E20" "  +3.86%  "
Set-TextValue "D21" "72.73"
Set-TextValue "E21" "  -0.73%  "
Set-TextValue "D22" "6.07"
Set-TextValue "E22" "  -1.22%  "
Set-TextValue "D23" "247.93"
Set-TextValue "E23" "  -0.91%  "
Set-TextValue "E25" "  -5.06%  "
Set-TextValue "E26" "  +21.25%  "
Set-TextValue "E27" "  -4.84%  "
Set-TextValue "D28" "9.81"
Set-TextValue "E28" "  -0.32%  "
Set-TextValue "D29" "173.15"
Set-TextValue "E29" "  +1.08%  "
Set-TextValue "D30" "21.00"
Set-TextValue "E30" "  +2.51%  "
Set-TextValue "E31" "  +1.04%  "
Set-TextValue "E32" "  -1.73%  "
Set-TextValue "D33" "0.125"
Set-TextValue "E33" "  +2.27%  "
Set-TextValue "E34" "  -2.55%  "
Set-TextValue "E35" "  +1.90%  "
Set-TextValue "D36" "4.94"
Set-TextValue "E36" "  -3.49%  "
Set-TextValue "E37" "  -0.72%  "
Set-TextValue "D38" "6.38"
Set-TextValue "E38" "  -3.46%  "
Set-TextValue "E39" "  -3.64%  "
Set-TextValue "E40" "  +2.51%  "
Set-TextValue "E41" "  +0.10%  "
Set-TextValue "D42" "8.62"
Set-TextValue "E42" "  +1.81%  "
Set-TextValue "D43" "0.000223"
Set-TextValue "E43" "  +0.84%  "
Set-TextValue "D44" "97.74"
Set-TextValue "E44" "  +0.10%  "
Set-TextValue "D45" "16.99"
Set-TextValue "E45" "  +0.96%  "
Set-TextValue "E46" "  -1.42%  "
Set-TextValue "D47" "0.0946"
Set-TextValue "E47" "  -1.84%  "
Set-TextValue "D48" "4.36"
Set-TextValue "E48" "  -6.71%  "
Set-TextValue "D49" "1.441.78"
Set-TextValue "E49" "  -2.76%  "
Set-TextValue "E50" "  -2.42%  "
Set-TextValue "E51" "  +1.50%  "
